$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (row 1, col A) ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 18:52"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 803018
$ws.Range("C4").Value = 10259
$ws.Range("D4").Value = 75317
$ws.Range("E4").Value = 684143
$ws.Range("G4").Value = 1044
$ws.Range("H4").Value = 43558

# --- Row 18: Suiza ---
$ws.Range("E18").Value = 7985
$ws.Range("G18").Value = 49
$ws.Range("H18").Value = 1478

# --- Row 22: Irlanda ---
$ws.Range("F22").Value = 152

# --- Rows 53-55: Sudafrica moves above Banglades & Egipto (sorted by Casos totales) ---
# Row 53 becomes Sudafrica with updated totals
$ws.Range("A53").Value = "Sudafrica"
$ws.Range("B53").Value = 3465
$ws.Range("C53").Value = 165
$ws.Range("D53").Value = 1055
$ws.Range("E53").Value = 2352
$ws.Range("F53").Value = 36
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 58

# Row 54 becomes Banglades (previous row 53 values)
$ws.Range("A54").Value = "Banglades"
$ws.Range("B54").Value = 3382
$ws.Range("C54").Value = 434
$ws.Range("D54").Value = 87
$ws.Range("E54").Value = 3185
$ws.Range("F54").Value = 1
$ws.Range("G54").Value = 9
$ws.Range("H54").Value = 110

# Row 55 becomes Egipto (previous row 54 values)
$ws.Range("A55").Value = "Egipto"
$ws.Range("B55").Value = 3333
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 821
$ws.Range("E55").Value = 2262
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 250

# --- Row 64: Kazajistan ---
$ws.Range("B64").Value = 1995
$ws.Range("C64").Value = 143
$ws.Range("E64").Value = 1487

# --- Row 105: Estado de Palestina ---
$ws.Range("B105").Value = 466
$ws.Range("C105").Value = 17
$ws.Range("E105").Value = 391

# --- Row 116: Isla de Man ---
$ws.Range("B116").Value = 307
$ws.Range("C116").Value = 7
$ws.Range("E116").Value = 98
